$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

$ws1.Range("H2").Value = 0.58841232396507
$ws1.Range("I2").Value = 0.280826277319723
$ws1.Range("O2").Value = 0.411640923743397

$ws1.Range("F3").Value = 0.603730324507612
$ws1.Range("G3").Value = 0.311559776320818

$ws1.Range("C4").Value = 0.637358737239101
$ws1.Range("D4").Value = 0.362731764845503
$ws1.Range("E4").Value = 1.0000905020846
$ws1.Range("J4").Value = 0.362698939830343
$ws1.Range("K4").Value = 0.311531582060349
$ws1.Range("L4").Value = 0.0153166143440778
$ws1.Range("M4").Value = 0.0489419839130541
$ws1.Range("N4").Value = 0.326848196404427

# --- Sheet 2: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

$ws2.Range("C2").Value = 0.362698939830343
$ws2.Range("D2").Value = 0.331220133480282
$ws2.Range("E2").Value = 0.394177746180405

$ws2.Range("C3").Value = 0.326848196404427
$ws2.Range("D3").Value = 0.297364592378915
$ws2.Range("E3").Value = 0.356331800429939

$ws2.Range("C4").Value = 0.411640923743397
$ws2.Range("D4").Value = 0.380403652512811
$ws2.Range("E4").Value = 0.442878194973983

$wb.Save()
